# Fix the leading-space typos in the header/text cells of the first sheet
# and update the active selection to B4, matching the corrected test asset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "trow"
$ws.Range("C1").Value = "frow"
$ws.Range("B2").Value = "Some text"
$ws.Range("B3").Value = "Some more test"

$null = $ws.Range("B4").Select()
